# Insert a new row at position 120 (shifts existing rows 120..141 down to 121..142)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(120).Insert()

# Fill in the new row 120 with the new weekly record
$ws.Range("A120").Value = 4
$ws.Range("B120").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C120").Value = "Los Lagos"
$ws.Range("D120").NumberFormat = $ws.Range("D121").NumberFormat
$ws.Range("D120").Value = 44474
$ws.Range("E120").Value = 10
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100108
$ws.Range("H120").Value = "Tropicales y subtropicales"
$ws.Range("I120").Value = 100108005
$ws.Range("J120").Value = "Piña"
$ws.Range("K120").Value = "Caramelo"
$ws.Range("L120").Value = "Primera"
$ws.Range("M120").Value = 80
$ws.Range("N120").Value = 22000
$ws.Range("O120").Value = 22000
$ws.Range("P120").Value = 22000
$ws.Range("Q120").Value = "$/caja 12 unidades"
$ws.Range("R120").Value = "Ecuador"
$ws.Range("S120").Value = 1833
$ws.Range("T120").Value = 12
